# ---------------------------------------------------------------------------
# training_plans.xlsx - apply the "Trucks with GPS" / "Employees equip with
# eletronic device" / "Deploy drones" training-plan sections (rows 10-23),
# rename sheet + reselect, per commit "wip plan de formation et plan mise ne
# place Leavitt Diamant".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- workbook-level tweaks --------------------------------------------------
$ws.Name = "Training plans"

# --- clear the old placeholder rows 10:13 (were just empty A10:A13 cells) --
$ws.Range("A10:J13").Clear()

# ===========================================================================
# Section 2 - "2 - Trucks with GPS"  (rows 10-14, header merged A10:A14)
# ===========================================================================
$sec2 = $ws.Range("A10:A14")
$sec2.Merge()
$sec2.Value = "2 - Trucks with GPS"
$sec2.WrapText = $true
$sec2.HorizontalAlignment = -4108   # xlCenter
$sec2.VerticalAlignment = -4108     # xlCenter
$sec2.Interior.Color = 13431551     # FFF2CC  (accent4 "gold", tint 0.8)
$ws.Rows.Item(10).RowHeight = 17
$ws.Rows.Item(11).RowHeight = 17
$ws.Rows.Item(12).RowHeight = 34
$ws.Rows.Item(13).RowHeight = 17
$ws.Rows.Item(14).RowHeight = 17

$ws.Range("B10").Value = "Team managers"
$ws.Range("C10").Value = "get truck location, how to identify a truck"
$ws.Range("D10").Value = "1h"
$ws.Range("E10").Value = "All team manager department"

$ws.Range("B11").Value = "Team manager referent"
$ws.Range("C11").Value = "advance gps settings, "
$ws.Range("D11").Value = "2h"
$ws.Range("E11").Value = "Team manager referent 1"

$ws.Range("B12").Value = "Truck's driver"
$ws.Range("C12").Value = "Check gps on, how to turn on off, how to repport a malfunction to the maintenance team"
$ws.Range("D12").Value = "1h"
$ws.Range("E12").Value = "All truck's driver"

$ws.Range("B13").Value = "Mainteance department"
$ws.Range("C13").Value = "How to install a gps on truck "
$ws.Range("D13").Value = "1h "
$ws.Range("E13").Value = "All maintenance department "

$ws.Range("B14").Value = "IT department "
$ws.Range("C14").Value = "Check truck's location "
$ws.Range("D14").Value = "1h "
$ws.Range("E14").Value = "All It department "

$dataB10 = $ws.Range("B10:B14,C10:C14,E10:E14")
$dataB10.WrapText = $true
$dataB10.VerticalAlignment = -4108

$dataD10 = $ws.Range("D10:D14")
$dataD10.WrapText = $true
$dataD10.VerticalAlignment = -4108
$dataD10.HorizontalAlignment = -4108

# ===========================================================================
# Section 3 - "3 - Employees equip with eletronic device" (rows 15-18)
# ===========================================================================
$sec3 = $ws.Range("A15:A18")
$sec3.Merge()
$sec3.Value = "3 - Employees equip with eletronic device"
$sec3.WrapText = $true
$sec3.HorizontalAlignment = -4108
$sec3.VerticalAlignment = -4108
$ws.Rows.Item(15).RowHeight = 34
$ws.Rows.Item(16).RowHeight = 17
$ws.Rows.Item(17).RowHeight = 34
$ws.Rows.Item(18).RowHeight = 17

# same fill as the existing "1 - use of the time manager application" header
# (A3:A9) -> copy/paste the format so it reuses that exact fill.
$ws.Range("A3").Copy() | Out-Null
$sec3.PasteSpecial(-4122) | Out-Null
$sec3.Value = "3 - Employees equip with eletronic device"

$ws.Range("B15").Value = "All employees"
$ws.Range("C15").Value = "Basic presentation "
$ws.Range("D15").Value = "1h"
$ws.Range("E15").Value = "All employees "

$merge1618 = $ws.Range("B16:B18")
$merge1618.Merge()
$merge1618.Value = "IT new recrue"
$merge1618.WrapText = $true
$merge1618.HorizontalAlignment = -4108
$merge1618.VerticalAlignment = -4108

$ws.Range("C16").Value = "Basic entreprise onboarding"
$ws.Range("D16").Value = "2h"
$ws.Range("E16").Value = "IT new recrue"

$ws.Range("C17").Value = "How execute preventive and currative maintenance plan "
$ws.Range("D17").Value = "3d"
$ws.Range("E17").Value = "IT new recrue"

$ws.Range("C18").Value = "Technical assistance posture"
$ws.Range("D18").Value = "1d"
$ws.Range("E18").Value = "IT new recrue"

$dataB15 = $ws.Range("B15,C15:C18,E15:E18")
$dataB15.WrapText = $true
$dataB15.VerticalAlignment = -4108

$dataD15 = $ws.Range("D15:D18")
$dataD15.WrapText = $true
$dataD15.VerticalAlignment = -4108
$dataD15.HorizontalAlignment = -4108

# ===========================================================================
# Section 4 - "4 - Deploy drones to monitor the street cleenliness" (19-23)
# ===========================================================================
$sec4 = $ws.Range("A19:A23")
$sec4.Merge()
$sec4.Value = "4 - Deploy drones to monitor the street cleenliness"
$sec4.WrapText = $true
$sec4.HorizontalAlignment = -4108
$sec4.VerticalAlignment = -4108
$sec4.Interior.Color = 16754931   # F3A8FF custom purple/pink
$ws.Rows.Item(19).RowHeight = 51
$ws.Rows.Item(20).RowHeight = 17
$ws.Rows.Item(21).RowHeight = 17
$ws.Rows.Item(22).RowHeight = 17
$ws.Rows.Item(23).RowHeight = 34

$ws.Range("B19").Value = "Drone pilot "
$ws.Range("C19").Value = "Basic entreprise onboarding"
$ws.Range("D19").Value = "2h"
$ws.Range("E19").Value = "2 pilots"

$ws.Range("B20").Value = "pilots + team manager"
$ws.Range("C20").Value = "how to define a good flight plan"
$ws.Range("D20").Value = "1d"
$ws.Range("E20").Value = "2 pilots and 1 team manager"

$ws.Range("B21").Value = "pilots + team manager"
$ws.Range("C21").Value = "report anomaly procedure "
$ws.Range("D21").Value = "1h"
$ws.Range("E21").Value = "2 pilots and 1 team manager"

$ws.Range("B22").Value = "maintenance Department"
$ws.Range("C22").Value = "Resolve basic failure on drones "
$ws.Range("D22").Value = "1d"
$ws.Range("E22").Value = "all maintenance department "

$ws.Range("B23").Value = "maintenance Department Referent"
$ws.Range("C23").Value = "Definition of a maintenance procedure with supplier / after sales service for serious breakdowns"
$ws.Range("D23").Value = "4h"
$ws.Range("E23").Value = "Maintenance departmenet referent"

# row 19 (B19,C19,E19) only wraps, no centering - matches the rest of the
# "row 20-23" body style otherwise used across the sheet
$row19cells = $ws.Range("B19,C19,E19")
$row19cells.WrapText = $true

$dataB20 = $ws.Range("B20:B23,C20:C23,E20:E23")
$dataB20.WrapText = $true
$dataB20.VerticalAlignment = -4108

$dataD19 = $ws.Range("D19:D23")
$dataD19.WrapText = $true
$dataD19.VerticalAlignment = -4108
$dataD19.HorizontalAlignment = -4108

# --- selection / scroll position matching the end state --------------------
$ws.Range("B19").Select()
$excel.ActiveWindow.ScrollRow = 6
